$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.504.93'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.11%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.098.79'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.45%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '384.19'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.25%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.86'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.20%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.541'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.76%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.02%  '

# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.88%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.89'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.37%  '

# Row 11
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.10%  '

# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.33%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.586.76'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +2.52%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.62'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.56%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.83'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.00%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.101.65'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +2.42%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '11.09'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +6.60%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.995'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.17%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '51.562.38'

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.35'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +9.65%  '

# Row 21
$ws.Range("B21").Value = 'InternetComputer(DFINITY)'
$ws.Range("C21").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.37'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.01%  '

# Row 22
$ws.Range("B22").Value = 'ShibaInu'
$ws.Range("C22").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0964'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.16%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.97'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.15%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '266.21'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.67%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.14'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.49%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.11'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.11%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27.05'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.31%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.24'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -3.44%  '

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.08%  '

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -2.94%  '

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.06%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.33'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.08%  '

# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.05%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0466'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.72%  '

# Row 35
$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.25'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.48%  '

# Row 36
$ws.Range("B36").Value = 'Toncoin'
$ws.Range("C36").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.04'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.59%  '

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.13%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.36'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.33%  '

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +4.23%  '

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.75%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '128.90'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.07%  '

# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -5.07%  '

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.37%  '

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.67%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '22.44'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.68%  '

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.69%  '

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +3.71%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.07'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.20%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.056.86'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.19%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0331'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +3.15%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.895'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +13.12%  '
